# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 14-15),
# pushing the existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 14.. down by two rows, duplicating the D2-style
# (date) formatting onto the freshly inserted rows.
$ws.Rows.Item(14).Resize(2).Insert()

# --- New row 14 : Primera, 50 bandejas @ 20000 --------------------------
$ws.Cells.Item(14,1).Value = 6
$ws.Cells.Item(14,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(14,3).Value = "Metropolitana"
$ws.Cells.Item(14,4).Value = 45037
$ws.Cells.Item(14,5).Value = 13
$ws.Cells.Item(14,6).Value = "Fruta"
$ws.Cells.Item(14,7).Value = 100101
$ws.Cells.Item(14,8).Value = "Berries"
$ws.Cells.Item(14,9).Value = 100101006
$ws.Cells.Item(14,10).Value = "Higo"
$ws.Cells.Item(14,11).Value = "Sin especificar"
$ws.Cells.Item(14,12).Value = "Primera"
$ws.Cells.Item(14,13).Value = 50
$ws.Cells.Item(14,14).Value = 20000
$ws.Cells.Item(14,15).Value = 20000
$ws.Cells.Item(14,16).Value = 20000
$ws.Cells.Item(14,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(14,18).Value = "Región Metropolitana"
$ws.Cells.Item(14,19).Value = 2857
$ws.Cells.Item(14,20).Value = 7

# --- New row 15 : Segunda, 50 bandejas @ 14000 ---------------------------
$ws.Cells.Item(15,1).Value = 6
$ws.Cells.Item(15,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(15,3).Value = "Metropolitana"
$ws.Cells.Item(15,4).Value = 45037
$ws.Cells.Item(15,5).Value = 13
$ws.Cells.Item(15,6).Value = "Fruta"
$ws.Cells.Item(15,7).Value = 100101
$ws.Cells.Item(15,8).Value = "Berries"
$ws.Cells.Item(15,9).Value = 100101006
$ws.Cells.Item(15,10).Value = "Higo"
$ws.Cells.Item(15,11).Value = "Sin especificar"
$ws.Cells.Item(15,12).Value = "Segunda"
$ws.Cells.Item(15,13).Value = 50
$ws.Cells.Item(15,14).Value = 14000
$ws.Cells.Item(15,15).Value = 14000
$ws.Cells.Item(15,16).Value = 14000
$ws.Cells.Item(15,17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(15,18).Value = "Región Metropolitana"
$ws.Cells.Item(15,19).Value = 2000
$ws.Cells.Item(15,20).Value = 7
